# Set the "treatment" value (column G) to "FNDR" for rows 2 through 255,
# which previously had no treatment recorded.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G255").Value = "FNDR"

# Reflect the active cell selection recorded in the saved workbook.
$ws.Range("G12").Select()
